$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(38, 1).Value = "SXT"
$ws.Cells.Item(38, 2).Value = "Year"
$ws.Cells.Item(38, 3).Value = [double]"8.904653999482134e-07"

$ws.Cells.Item(39, 1).Value = "SXT"
$ws.Cells.Item(39, 2).Value = "Specimen_type"
$ws.Cells.Item(39, 3).Value = [double]"7.953466886058472e-06"

$ws.Cells.Item(40, 1).Value = "SXT"
$ws.Cells.Item(40, 2).Value = "Gender"
$ws.Cells.Item(40, 3).Value = [double]"5.429345858914697e-11"

$ws.Cells.Item(41, 1).Value = "SXT"
$ws.Cells.Item(41, 2).Value = "Hospital:Ward_ED_ICU"
$ws.Cells.Item(41, 3).Value = [double]"1.862655002196334e-08"
